# ---------------------------------------------------------------------------
# Markov quick example -> new input layout ("ScenarioA") + DB-style metadata
# columns (id / description / unit / data package / data source), matching
# the commit "Adjust markov quick example to new input / Add option to
# define output-file path for unit-commitment results".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the worksheet & fix up the dependent defined name ----------
$ws.Name = "ScenarioA"

foreach ($n in $wb.Names) {
    if ($n.Name -eq "weight_rp") {
        $n.RefersTo = "=ScenarioA!`$B`$5:`$D`$15"
    }
}

# --- 2. Drop the old helper comments that sat on cells which are being ----
#        repurposed/removed so they don't linger on the wrong cell.
$ws.Range("C4").Comment.Delete()
$ws.Range("C5").Comment.Delete()

# --- 3. Title & format rows -------------------------------------------
$ws.Range("B1").Value = "Power - Weights Representative Periods"

$ws.Range("B2").Value = "Format:"
$ws.Range("C2").Value = "v0.0.2"

# --- 4. New 5-column metadata header block (rows 3-7) ----------------------
$ws.Range("B3").Value = "Database ID"
$ws.Range("C3").Value = "rp"
$ws.Range("D3").Value = "representatives periods weight"
$ws.Range("E3").Value = "Data Package"
$ws.Range("F3").Value = "Data Source"

$ws.Range("B4").Value = "id"
$ws.Range("C4").Value = "rp"
$ws.Range("D4").Value = "pWeight_rp"
$ws.Range("E4").Value = "dataPackage"
$ws.Range("F4").Value = "dataSource"

$ws.Range("B5").Value = "ID within database"
$ws.Range("C5").Value = "Representative period"
$ws.Range("D5").Value = "Weight of representative period within full timespan"
$ws.Range("E5").Value = "Which package this belongs to"
$ws.Range("F5").Value = "Where the data for the entry comes from"

$ws.Range("B6").Value = "Filled automatically by database"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "Scenario-dependent"
$ws.Range("E6").Value = "Scenario-dependent"
$ws.Range("F6").Value = "Scenario-dependent"

$ws.Range("B7").Value = "[db-key]"
$ws.Range("C7").Value = "[rp]"
$ws.Range("D7").Value = "[h]"
$ws.Range("E7").Value = "[DataPackage]"
$ws.Range("F7").Value = "[DataSource]"

# --- 5. Data rows (8-13): rp01..rp06, weight 1, test data package/source --
$rpNames = @("rp01", "rp02", "rp03", "rp04", "rp05", "rp06")
for ($i = 0; $i -lt $rpNames.Length; $i++) {
    $r = 8 + $i
    $ws.Cells.Item($r, 3).Value = $rpNames[$i]
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = "TestPackage1"
    $ws.Cells.Item($r, 6).Value = "TestSource1"
}

# --- 6. Comments describing each row label (B3:B7) -------------------------
$ws.Range("B3").Comment.Text("Readable Name")
$ws.Range("B4").AddComment("Value specifier in DB")
$ws.Range("B5").AddComment("Description")
$ws.Range("B6").AddComment("Details on database behavior")
$ws.Range("B7").AddComment("Unit or valid values")

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# Column widths
$ws.Columns("B").ColumnWidth = 16.5
$ws.Columns("C").ColumnWidth = 21.5
$ws.Range("D1:F1").EntireColumn.ColumnWidth = 31.666666666666668

# Title row (B1): bold white Aptos 18 on teal fill, vertically centered
$titleRange = $ws.Range("B1")
$titleRange.Font.Name = "Aptos"
$titleRange.Font.Size = 18
$titleRange.Font.Bold = $true
$titleRange.Font.Color = 16777215
$titleRange.Interior.Color = 8421376
$titleRange.VerticalAlignment = -4108
$ws.Rows(1).RowHeight = 30

# "Format:" label (B2): italic Aptos 11, right aligned
$fmtLabel = $ws.Range("B2")
$fmtLabel.Font.Name = "Aptos"
$fmtLabel.Font.Size = 11
$fmtLabel.Font.Italic = $true
$fmtLabel.HorizontalAlignment = -4152

# "v0.0.2" value (C2): italic Aptos 11
$fmtValue = $ws.Range("C2")
$fmtValue.Font.Name = "Aptos"
$fmtValue.Font.Size = 11
$fmtValue.Font.Italic = $true

# Row 3 header (Database ID / rp / ... ): bold Aptos 11 on light blue fill
$row3 = $ws.Range("B3:F3")
$row3.Font.Name = "Aptos"
$row3.Font.Size = 11
$row3.Font.Bold = $true
$row3.Interior.Color = 16050907
$row3.HorizontalAlignment = -4131

# Row 4 header (id / rp / ...): bold Aptos 11 on light grey fill
$row4 = $ws.Range("B4:F4")
$row4.Font.Name = "Aptos"
$row4.Font.Size = 11
$row4.Font.Bold = $true
$row4.Interior.Color = 15921906
$row4.HorizontalAlignment = -4131

# Row 5 (descriptions): italic Aptos 11 on light fill, wrap, top aligned
$row5 = $ws.Range("B5:F5")
$row5.Font.Name = "Aptos"
$row5.Font.Size = 11
$row5.Font.Italic = $true
$row5.Interior.Color = 15921906
$row5.HorizontalAlignment = -4131
$row5.VerticalAlignment = -4160
$row5.WrapText = $true
$ws.Rows(5).RowHeight = 30

# Row 6 (database behavior / scenario-dependent): italic Aptos 11, grey fill, wrap
$row6 = $ws.Range("B6:F6")
$row6.Font.Name = "Aptos"
$row6.Font.Size = 11
$row6.Font.Italic = $true
$row6.Interior.Color = 14211288
$row6.HorizontalAlignment = -4131
$row6.VerticalAlignment = -4160
$row6.WrapText = $true
$ws.Rows(6).RowHeight = 45

# Row 7 ([db-key] / [rp] / [h] / ...): italic Aptos 11, light fill, wrap, indent
$row7 = $ws.Range("B7:F7")
$row7.Font.Name = "Aptos"
$row7.Font.Size = 11
$row7.Font.Italic = $true
$row7.Interior.Color = 15921906
$row7.HorizontalAlignment = -4131
$row7.VerticalAlignment = -4160
$row7.WrapText = $true
$row7.IndentLevel = 1

# B8:B13 left blank but keep consistent plain styling
$colB813 = $ws.Range("B8:B13")
$colB813.Font.Name = "Aptos"
$colB813.Font.Size = 11
$colB813.Font.Italic = $true
$colB813.Interior.Color = 15921906
$colB813.HorizontalAlignment = -4131
$colB813.VerticalAlignment = -4160
$colB813.WrapText = $true
$colB813.IndentLevel = 1

# C8:C13 (rp01..rp06): plain Aptos 11 on green fill, indented
$colC813 = $ws.Range("C8:C13")
$colC813.Font.Name = "Aptos"
$colC813.Font.Size = 11
$colC813.Interior.Color = 10284031
$colC813.HorizontalAlignment = -4131
$colC813.IndentLevel = 1

# D8:D13 (weight = 1): plain Aptos 11 on light blue fill, right aligned, integer format
$colD813 = $ws.Range("D8:D13")
$colD813.Font.Name = "Aptos"
$colD813.Font.Size = 11
$colD813.Interior.Color = 15129530
$colD813.HorizontalAlignment = -4152
$colD813.NumberFormat = "0"

# E8:F13 (TestPackage1 / TestSource1): plain Aptos 11 on light blue fill, indented
$colEF813 = $ws.Range("E8:F13")
$colEF813.Font.Name = "Aptos"
$colEF813.Font.Size = 11
$colEF813.Interior.Color = 15129530
$colEF813.HorizontalAlignment = -4131
$colEF813.IndentLevel = 1

# --- 7. Freeze panes below the header block, to the right of column B -----
$ws.Range("C8").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 100

$ws.Range("A1").Select()
